$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsData = @"
1847,461.2497333333332,441.2768041237113,14.24333333333333
1848,461.4997333333332,310.5088917525773,14.55666666666667
1849,461.7498000000001,193.4444587628866,14.41466666666667
1850,461.9997333333332,85.19974226804123,14.42533333333333
1851,462.2497333333334,20.58994845360824,14.46366666666666
1852,462.4997333333334,3.644072164948454,14.35733333333334
1853,462.7497333333333,2.899871134020619,14.21966666666666
1854,462.9997333333333,3.404510309278351,14.08499999999999
1855,463.2498166666665,3.66971649484536,14.05533333333333
1856,463.4997333333333,3.319072164948453,14.219
1857,463.7497333333333,3.147938144329897,14.35033333333333
1858,463.9997999999999,2.822938144329897,14.255
1859,464.2497333333333,2.532087628865979,14.37533333333333
1860,464.4998166666667,2.215592783505155,14.13933333333334
1861,464.7497333333332,2.130025773195876,14.31233333333334
1862,465.00395,2.232603092783505,14.52433333333334
1863,465.2524333333332,2.12139175257732,14.191
1864,465.5009499999999,1.804896907216495,14.263
1865,465.7497333333332,1.479896907216495,14.42966666666667
1866,465.9997833333333,1.197551546391752,14.429
1867,466.2497333333334,1.308762886597938,13.80866666666666
1868,466.4998,1.120618556701031,14.20166666666666
1869,466.7497333333333,0.940979381443299,14.06233333333334
1870,466.9997499999999,0.6586340206185567,13.628
1871,467.2498166666667,0.581701030927835,13.19466666666667
1872,467.4997333333333,0.5474226804123711,13.017
1873,467.7497999999999,0.8041237113402062,12.89333333333333
1874,467.9997333333333,0.8297680412371133,13.003
1875,468.2497333333333,0.7185567010309278,13.02866666666666
1876,468.4997333333332,0.5132731958762886,13.14033333333333
1877,468.7497333333332,0.3677835051546391,13.37733333333334
1878,468.9997999999998,0.2480670103092784,13.49066666666667
1879,469.2497666666665,0.1453608247422681,13.44566666666666
1880,469.4997333333332,0.2565721649484536,13.609
1881,469.7498166666666,0.1967783505154639,13.64366666666667
1882,469.9997333333334,,13.629
1883,470.2497333333333,,13.83133333333334
1884,470.4997333333333,,13.77866666666667
1885,470.7497333333333,,13.31633333333334
1886,470.9997499999998,,12.89766666666667
1887,471.2497333333333,,12.24466666666667
1888,471.4997333333333,,11.46333333333333
1889,471.7497666666666,,10.33066666666667
1890,471.9997333333333,,9.125666666666671
1891,472.24975,,8.065333333333328
1892,472.4997833333333,,7.063333333333333
1893,472.7497333333332,,5.945999999999994
1894,472.9998000000001,,4.903666666666677
1895,473.2497333333332,,3.625333333333334
1896,473.4997333333334,,2.788666666666668
1897,473.7497333333334,,2.193000000000001
1898,473.9997333333333,,2.002333333333333
1899,474.2497833333332,,1.773666666666667
1900,474.4997333333333,,1.671666666666663
1901,474.7498166666667,,1.515000000000001
1902,474.9997333333333,,1.391999999999999
1903,475.2497333333333,,1.257999999999999
1904,475.4997333333333,,1.094999999999999
1905,475.7498166666667,,0.9669999999999952
1906,475.9997333333332,,0.7516666666666652
1907,476.24975,,0.6396666666666633
1908,476.4997333333332,,0.4913333333333334
1909,476.7498000000001,,0.416999999999998
1910,476.9997333333332,,0.3343333333333334
1911,477.2497833333333,,0.2763333333333371
1912,477.4997333333334,,0.2199999999999989
1913,477.7498,,0.1373333333333342
1914,477.9997333333333,,0.02799999999999869
1915,478.2498166666665,,0.0313333333333361
1916,478.4997333333333,,-0.02033333333333331
1917,478.7497833333334,,-0.1350000000000016
1918,478.9997333333333,,-0.2119999999999997
1919,479.2497333333333,,-0.1906666666666652
1920,479.4997333333333,,-0.1506666666666661
1921,479.7497833333333,,-0.1333333333333293
1922,479.9997333333332,,-0.2056666666666658
1923,480.2497333333332,,-0.251333333333335
1924,480.4998000000001,,-0.2966666666666633
1925,480.7497333333332,,-0.1990000000000016
1926,480.9997833333333,,-0.3256666666666668
1927,481.2497333333334,,-0.3006666666666682
1928,481.4997666666667,,-0.396333333333331
1929,481.7497333333333,,-0.4209999999999994
1930,481.9997333333333,,-0.320333333333334
1931,482.2498166666667,,-0.3866666666666667
1932,482.4997333333333,,-0.3989999999999974
1933,482.7497333333333,,-0.3926666666666669
1934,482.9997333333333,,-0.3466666666666676
1935,483.2497833333333,,-0.378666666666664
1936,483.4998166666667,,-0.4083333333333314
1937,483.7497333333332,,-0.4260000000000055
1938,483.9997666666665,,-0.3840000000000003
1939,484.2497333333332,,-0.3853333333333353
1940,484.4997499999999,,-0.4140000000000015
1941,484.7497833333333,,-0.3880000000000017
1942,484.9997333333334,,-0.3860000000000028
1943,485.2497666666667,,-0.4506666666666668
1944,485.4997333333333,,-0.4106666666666676
1945,485.7497499999999,,-0.3816666666666677
1946,485.9997333333333,,-0.3866666666666667
1947,486.2497499999998,,-0.3153333333333386
1948,486.4997333333333,,-0.277666666666665
1949,486.7497999999999,,-0.3026666666666671
1950,486.9997333333333,,-0.277000000000001
1951,487.2498166666667,,-0.277000000000001
1952,487.4997333333332,,-0.1986666666666643
1953,487.7497999999998,,-0.2473333333333301
1954,487.9998000000001,,-0.1993333333333354
1955,488.2498166666666,,-0.2013333333333343
1956,488.4997333333334,,-0.2360000000000007
1957,488.7497333333334,,-0.2070000000000007
1958,488.9997333333333,,-0.1166666666666671
"@

$lines = $rowsData -split "`r?`n"

foreach ($line in $lines) {
    if ($line.Trim() -eq "") { continue }
    $parts = $line -split ","
    $r = [int]$parts[0]
    $aVal = [double]$parts[1]
    $bText = $parts[2]
    $cVal = [double]$parts[3]

    $ws.Cells.Item($r, 1).Value = $aVal
    if ($bText -ne "") {
        $ws.Cells.Item($r, 2).Value = [double]$bText
    }
    $ws.Cells.Item($r, 3).Value = $cVal
}

# Apply the same style as the rest of column A (bold font, thin border, center/top alignment)
$styleRange = $ws.Range("A1847:A1958")
$styleRange.Font.Bold = $true
$styleRange.Borders.LineStyle = 1
$styleRange.HorizontalAlignment = -4108
$styleRange.VerticalAlignment = -4160

Write-Host "Added rows 1847-1958"
